# Applies the "Book insert query added" commit to scrum-planering.xlsx
$wb = $excel.ActiveWorkbook

# --- Product Backlog sheet ---
$ws1 = $wb.Worksheets.Item("Product Backlog")
$ws1.Range("C10").Value = 1
$ws1.Range("C11").Value = 1

# --- Sprint 2 sheet ("Sprint 2 xx.xx - xx.xx") ---
$ws3 = $wb.Worksheets.Item("Sprint 2 xx.xx - xx.xx")
$ws3.Range("E13").Value = 360
$ws3.Range("G13").Value = "Rasmus"
$ws3.Range("H13").Value = (Get-Date -Year 2024 -Month 11 -Day 19 -Hour 0 -Minute 0 -Second 0).Date
$ws3.Range("E14").Value = 60
$ws3.Range("G14").Value = "Rasmus"
$ws3.Range("H14").Value = (Get-Date -Year 2024 -Month 11 -Day 20 -Hour 0 -Minute 0 -Second 0).Date
$ws3.Range("E12").Formula = "=SUM(E13:E17)/60"
$ws3.Range("D18").Formula = "=SUM(D6+D12)"
$ws3.Range("E18").Formula = "=SUM(E6+E12)"

# --- Sprint 3 sheet ("Sprint 3 xx.xx - xx.xx") ---
$ws4 = $wb.Worksheets.Item("Sprint 3 xx.xx - xx.xx")
$ws4.Range("D13").Value = 180
$ws4.Range("D17").Formula = "=SUM(D6+D12)"
$ws4.Range("D19").Value = 2

# --- Sprint 4 sheet ("Sprint 4 xx.xx - xx.xx") ---
$ws5 = $wb.Worksheets.Item("Sprint 4 xx.xx - xx.xx")
$ws5.Range("I7").Value = "Prepared querys, htacces, block user access, sanitize input and ?"

$wb.Save()
